$wb = $excel.ActiveWorkbook

# --- storage_data sheet: fix energy balance values ---
$wsStorage = $wb.Worksheets.Item("storage_data")
# initially_charged (M2): was text "false", now boolean TRUE
$wsStorage.Range("M2").Value = $true
# init_charging_temp (N2): 50 -> 35
$wsStorage.Range("N2").Value = 35

# --- house_data sheet: fix autarky_thermal value ---
$wsHouse = $wb.Worksheets.Item("house_data")
$wsHouse.Range("G2").Value = 0.16

# --- mediator_data sheet: add repetition_period column ---
$wsMediator = $wb.Worksheets.Item("mediator_data")
$wsMediator.Range("I1").Value = "repetition_period"
$wsMediator.Range("I1").HorizontalAlignment = 1
$wsMediator.Range("I2").Value = 3
$wsMediator.Columns.Item(9).AutoFit()

# --- restore selection / active sheet state ---
[void]$wsHouse.Range("G9").Select()
[void]$wsMediator.Range("E9").Select()

# storage_data becomes the active/selected tab with N9 selected
[void]$wsStorage.Activate()
[void]$wsStorage.Range("N9").Select()
